# Update "想去人数" (F column) values on sheet "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9619
$ws1.Range("F3").Value = 213
$ws1.Range("F4").Value = 31
$ws1.Range("F5").Value = 542
$ws1.Range("F6").Value = 466

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9619
$ws4.Range("F3").Value = 213
$ws4.Range("F4").Value = 31
$ws4.Range("F5").Value = 542
$ws4.Range("F7").Value = 466
